$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host ("Sheet name: " + $ws.Name)
$v = $ws.Range("A1").Value
Write-Host ("A1 type: " + $v.GetType())
Write-Host ("A1: $v")
Write-Host ("I1: " + $ws.Range("I1").Value)
Write-Host ("Cells I1: " + $ws.Cells.Item(1,9).Value)
